$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous data block (rows 2-13) so stale rows/strings don't linger.
$ws.Range("A2:T13").Clear()

$data = @(
    @("ECs","Trf","Tfrc","ECs",2,0.6666666666666666,0.4549493333333334,1.364848,0.004035600264670243,0.004035600264670243,3,1,11.68484733333333,35.054542,0.266948203543434,0.266948203543434,5.316013504401778,47.844121539616,0.001077296240873128,0.001077296240873128),
    @("ECs","Trf","Tfrc","FAPs",2,0.6666666666666666,0.4549493333333334,1.364848,0.004035600264670243,0.004035600264670243,3,1,4.833228666666667,14.499686,0.1104183626088705,0.1104183626088705,2.198874159747556,19.789867437728,0.0004456043733688127,0.0004456043733688126),
    @("ECs","Trf","Tfrc","M2",2,0.6666666666666666,0.4549493333333334,1.364848,0.004035600264670243,0.004035600264670243,3,1,14.330788,42.99236399999999,0.3273964993148507,0.3273964993148507,6.519782446741333,58.678042020672,0.001321241399287123,0.001321241399287122),
    @("ECs","Trf","Tfrc","sCs",2,0.6666666666666666,0.4549493333333334,1.364848,0.004035600264670243,0.004035600264670243,3,1,12.92310066666666,38.769302,0.2952369345328449,0.2952369345328449,5.879356032899556,52.914204296096,0.00119145825114118,0.00119145825114118),
    @("FAPs","Trf","Tfrc","ECs",3,1,19.60581866666667,58.81745600000001,0.1739122165990868,0.1739122165990868,3,1,11.68484733333333,35.054542,0.266948203543434,0.266948203543434,229.0909979650169,2061.818981685152,0.04642555379538279,0.04642555379538279),
    @("FAPs","Trf","Tfrc","FAPs",3,1,19.60581866666667,58.81745600000001,0.1739122165990868,0.1739122165990868,3,1,4.833228666666667,14.499686,0.1104183626088705,0.1104183626088705,94.75940481320178,852.8346433188161,0.01920310219455039,0.01920310219455039),
    @("FAPs","Trf","Tfrc","M2",3,1,19.60581866666667,58.81745600000001,0.1739122165990868,0.1739122165990868,3,1,14.330788,42.99236399999999,0.3273964993148507,0.3273964993148507,280.9668308784426,2528.701477905984,0.05693825090262707,0.05693825090262707),
    @("FAPs","Trf","Tfrc","sCs",3,1,19.60581866666667,58.81745600000001,0.1739122165990868,0.1739122165990868,3,1,12.92310066666666,38.769302,0.2952369345328449,0.2952369345328449,253.3679682817458,2280.311714535712,0.05134530970652652,0.05134530970652652),
    @("M2","Trf","Tfrc","ECs",3,1,91.29211166666666,273.876335,0.8098010987398714,0.8098010987398714,3,1,11.68484733333333,35.054542,0.266948203543434,0.266948203543434,1066.734387562619,9600.609488063568,0.2161749485361076,0.2161749485361076),
    @("M2","Trf","Tfrc","FAPs",3,1,91.29211166666666,273.876335,0.8098010987398714,0.8098010987398714,3,1,4.833228666666667,14.499686,0.1104183626088705,0.1104183626088705,441.2356511478678,3971.12086033081,0.08941691136172086,0.08941691136172086),
    @("M2","Trf","Tfrc","M2",3,1,91.29211166666666,273.876335,0.8098010987398714,0.8098010987398714,3,1,14.330788,42.99236399999999,0.3273964993148507,0.3273964993148507,1308.287898367326,11774.59108530594,0.2651260448687536,0.2651260448687536),
    @("M2","Trf","Tfrc","sCs",3,1,91.29211166666666,273.876335,0.8098010987398714,0.8098010987398714,3,1,12.92310066666666,38.769302,0.2952369345328449,0.2952369345328449,1179.777149140908,10617.99434226817,0.2390831939732893,0.2390831939732893),
    @("sCs","Trf","Tfrc","ECs",3,1,1.381113666666667,4.143340999999999,0.01225108439637166,0.01225108439637166,3,1,11.68484733333333,35.054542,0.266948203543434,0.266948203543434,16.13810234498022,145.242921104822,0.00327040497107041,0.003270404971070409),
    @("sCs","Trf","Tfrc","FAPs",3,1,1.381113666666667,4.143340999999999,0.01225108439637166,0.01225108439637166,3,1,4.833228666666667,14.499686,0.1104183626088705,0.1104183626088705,6.675238165658445,60.077143490926,0.001352744679230441,0.001352744679230441),
    @("sCs","Trf","Tfrc","M2",3,1,1.381113666666667,4.143340999999999,0.01225108439637166,0.01225108439637166,3,1,14.330788,42.99236399999999,0.3273964993148507,0.3273964993148507,19.79244716090266,178.132024448124,0.004010962144182872,0.004010962144182871),
    @("sCs","Trf","Tfrc","sCs",3,1,1.381113666666667,4.143340999999999,0.01225108439637166,0.01225108439637166,3,1,12.92310066666666,38.769302,0.2952369345328449,0.2952369345328449,17.84827094644244,160.634438517982,0.003616972601887937,0.003616972601887937)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
